$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 42632.87872685185
$ws.Range("B3").Value = -1
$ws.Range("C3").Value = "Neutral"
$ws.Range("D3").Value = 22
$ws.Range("E3").Value = 16754
$ws.Range("F3").Value = 2640
$ws.Range("G3").Value = 59
$ws.Range("H3").Value = 36
$ws.Range("I3").Value = 92
$ws.Range("J3").Value = 7
$ws.Range("K3").Value = 18455
$ws.Range("L3").Value = 371
$ws.Range("M3").Value = 225
$ws.Range("N3").Value = 49
$ws.Range("O3").Value = 4
$ws.Range("P3").Value = "Noun"
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = -31.57
$ws.Range("S3").Value = -0.0872
$ws.Range("T3").Value = -0.74
$ws.Range("U3").Value = 6.75
$ws.Range("V3").Value = 1.88
$ws.Range("W3").Value = 0

$ws.Range("A3").NumberFormat = "m/d/yy h:mm"
$ws.Range("S3").NumberFormat = "0.00%"
